# SEM Update on August 18
#
# Appends a new month column (07/01/2023) to both TABLE_1 (levels) and
# TABLE_2 (year-over-year pct change), and revises a handful of the most
# recent months' figures to match the refreshed source release.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("TABLE_1")
$ws2 = $wb.Worksheets.Item("TABLE_2")

# --- TABLE_1: new header cell (EJ4) ---
$c = $ws1.Cells.Item(4, 140)
$c.NumberFormat = "@"
$c.Value2 = "07/01/2023"
$c.Style = "Normal"

# --- TABLE_2: new header cell (DX4) ---
$c = $ws2.Cells.Item(4, 128)
$c.NumberFormat = "@"
$c.Value2 = "07/01/2023"
$c.Style = "Normal"

# --- TABLE_1: revisions to existing months + new EJ column ---
$ws1.Cells.Item(5, 138).Value2 = 10720.3
$ws1.Cells.Item(5, 139).Value2 = 10139.8
$ws1.Cells.Item(5, 140).Value2 = 9046.6
$ws1.Cells.Item(6, 139).Value2 = 169.2
$ws1.Cells.Item(6, 140).Value2 = 160.8
$ws1.Cells.Item(7, 139).Value2 = 22.9
$ws1.Cells.Item(7, 140).Value2 = 19.3
$ws1.Cells.Item(8, 139).Value2 = 152.7
$ws1.Cells.Item(8, 140).Value2 = 151.9
$ws1.Cells.Item(9, 140).Value2 = 82.2
$ws1.Cells.Item(10, 139).Value2 = 1234.6
$ws1.Cells.Item(10, 140).Value2 = 1074.5
$ws1.Cells.Item(11, 139).Value2 = 211.2
$ws1.Cells.Item(11, 140).Value2 = 197
$ws1.Cells.Item(12, 139).Value2 = 114.2
$ws1.Cells.Item(12, 140).Value2 = 96.7
$ws1.Cells.Item(13, 139).Value2 = 35.2
$ws1.Cells.Item(13, 140).Value2 = 34.1
$ws1.Cells.Item(15, 140).Value2 = 393
$ws1.Cells.Item(16, 139).Value2 = 335.3
$ws1.Cells.Item(16, 140).Value2 = 320
$ws1.Cells.Item(18, 139).Value2 = 61.4
$ws1.Cells.Item(18, 140).Value2 = 54.9
$ws1.Cells.Item(19, 139).Value2 = 408.7
$ws1.Cells.Item(19, 140).Value2 = 380.6
$ws1.Cells.Item(20, 139).Value2 = 185.5
$ws1.Cells.Item(20, 140).Value2 = 167.8
$ws1.Cells.Item(21, 139).Value2 = 132.2
$ws1.Cells.Item(21, 140).Value2 = 114.6
$ws1.Cells.Item(22, 139).Value2 = 119.9
$ws1.Cells.Item(22, 140).Value2 = 102.6
$ws1.Cells.Item(23, 140).Value2 = 126.9
$ws1.Cells.Item(24, 139).Value2 = 139.1
$ws1.Cells.Item(24, 140).Value2 = 131.9
$ws1.Cells.Item(25, 139).Value2 = 47.7
$ws1.Cells.Item(25, 140).Value2 = 40.7
$ws1.Cells.Item(26, 140).Value2 = 196.1
$ws1.Cells.Item(27, 139).Value2 = 239.7
$ws1.Cells.Item(27, 140).Value2 = 212.6
$ws1.Cells.Item(28, 140).Value2 = 263.6
$ws1.Cells.Item(29, 139).Value2 = 197.7
$ws1.Cells.Item(29, 140).Value2 = 167.5
$ws1.Cells.Item(30, 139).Value2 = 98.7
$ws1.Cells.Item(30, 140).Value2 = 93.3
$ws1.Cells.Item(32, 140).Value2 = 32.1
$ws1.Cells.Item(33, 139).Value2 = 90.5
$ws1.Cells.Item(33, 140).Value2 = 78.7
$ws1.Cells.Item(34, 140).Value2 = 70.6
$ws1.Cells.Item(35, 139).Value2 = 39.1
$ws1.Cells.Item(35, 140).Value2 = 32.3
$ws1.Cells.Item(36, 139).Value2 = 303.6
$ws1.Cells.Item(36, 140).Value2 = 257.1
$ws1.Cells.Item(37, 139).Value2 = 71
$ws1.Cells.Item(37, 140).Value2 = 65.3
$ws1.Cells.Item(38, 139).Value2 = 659.8
$ws1.Cells.Item(38, 140).Value2 = 538.4
$ws1.Cells.Item(39, 140).Value2 = 237.3
$ws1.Cells.Item(40, 139).Value2 = 32.4
$ws1.Cells.Item(40, 140).Value2 = 26.1
$ws1.Cells.Item(41, 139).Value2 = 357.7
$ws1.Cells.Item(41, 140).Value2 = 348.9
$ws1.Cells.Item(42, 139).Value2 = 134.8
$ws1.Cells.Item(42, 140).Value2 = 125.6
$ws1.Cells.Item(43, 139).Value2 = 141.6
$ws1.Cells.Item(43, 140).Value2 = 111.3
$ws1.Cells.Item(44, 139).Value2 = 302.2
$ws1.Cells.Item(44, 140).Value2 = 271.3
$ws1.Cells.Item(45, 139).Value2 = 27.3
$ws1.Cells.Item(45, 140).Value2 = 22.9
$ws1.Cells.Item(46, 139).Value2 = 155.8
$ws1.Cells.Item(46, 140).Value2 = 147.1
$ws1.Cells.Item(47, 139).Value2 = 34
$ws1.Cells.Item(47, 140).Value2 = 28.4
$ws1.Cells.Item(48, 139).Value2 = 192.8
$ws1.Cells.Item(48, 140).Value2 = 168.2
$ws1.Cells.Item(49, 139).Value2 = 1133.6
$ws1.Cells.Item(49, 140).Value2 = 1061
$ws1.Cells.Item(50, 140).Value2 = 107.5
$ws1.Cells.Item(51, 139).Value2 = 29.4
$ws1.Cells.Item(51, 140).Value2 = 26.1
$ws1.Cells.Item(52, 139).Value2 = 305.8
$ws1.Cells.Item(52, 140).Value2 = 273.2
$ws1.Cells.Item(53, 139).Value2 = 249.5
$ws1.Cells.Item(53, 140).Value2 = 230.2
$ws1.Cells.Item(54, 139).Value2 = 55.9
$ws1.Cells.Item(54, 140).Value2 = 51.7
$ws1.Cells.Item(55, 139).Value2 = 189.7
$ws1.Cells.Item(55, 140).Value2 = 168.1
$ws1.Cells.Item(56, 139).Value2 = 28.1
$ws1.Cells.Item(56, 140).Value2 = 24.1

# --- TABLE_2: revisions to existing months + new DX column ---
$ws2.Cells.Item(5, 126).Value2 = 2.36326483843862
$ws2.Cells.Item(5, 127).Value2 = 2.67112191170514
$ws2.Cells.Item(5, 128).Value2 = 1.785573645068
$ws2.Cells.Item(6, 127).Value2 = 0.594530321046373
$ws2.Cells.Item(6, 128).Value2 = 0.06222775357811
$ws2.Cells.Item(7, 127).Value2 = 4.09090909090908
$ws2.Cells.Item(7, 128).Value2 = 2.11640211640213
$ws2.Cells.Item(8, 127).Value2 = 1.73217854763491
$ws2.Cells.Item(8, 128).Value2 = 3.47411444141689
$ws2.Cells.Item(9, 128).Value2 = 0.611995104039185
$ws2.Cells.Item(10, 127).Value2 = 2.72069223729094
$ws2.Cells.Item(10, 128).Value2 = 3.08932169241102
$ws2.Cells.Item(11, 127).Value2 = 3.32681017612524
$ws2.Cells.Item(11, 128).Value2 = 5.74342458400429
$ws2.Cells.Item(12, 127).Value2 = 2.42152466367713
$ws2.Cells.Item(12, 128).Value2 = -0.309278350515461
$ws2.Cells.Item(13, 127).Value2 = 0.57142857142858
$ws2.Cells.Item(13, 128).Value2 = 0.887573964497054
$ws2.Cells.Item(15, 128).Value2 = 1.83985488468516
$ws2.Cells.Item(16, 127).Value2 = 3.51960481630134
$ws2.Cells.Item(16, 128).Value2 = 3.55987055016181
$ws2.Cells.Item(18, 127).Value2 = 5.67986230636834
$ws2.Cells.Item(18, 128).Value2 = 6.60194174757281
$ws2.Cells.Item(19, 127).Value2 = 3.65204159269591
$ws2.Cells.Item(19, 128).Value2 = 3.42391304347827
$ws2.Cells.Item(20, 127).Value2 = -2.87958115183246
$ws2.Cells.Item(20, 128).Value2 = -1.46799765120376
$ws2.Cells.Item(21, 127).Value2 = 0.379650721336371
$ws2.Cells.Item(21, 128).Value2 = 0
$ws2.Cells.Item(22, 127).Value2 = -0.580431177446093
$ws2.Cells.Item(22, 128).Value2 = 0.686947988223738
$ws2.Cells.Item(23, 128).Value2 = 3.42298288508559
$ws2.Cells.Item(24, 127).Value2 = 3.34323922734027
$ws2.Cells.Item(24, 128).Value2 = 2.08978328173376
$ws2.Cells.Item(25, 127).Value2 = 0.632911392405072
$ws2.Cells.Item(25, 128).Value2 = 0.742574257425753
$ws2.Cells.Item(26, 128).Value2 = 2.56276150627615
$ws2.Cells.Item(27, 127).Value2 = 2.43589743589743
$ws2.Cells.Item(27, 128).Value2 = 2.70531400966183
$ws2.Cells.Item(28, 128).Value2 = 5.73606097071801
$ws2.Cells.Item(29, 127).Value2 = 1.07361963190182
$ws2.Cells.Item(29, 128).Value2 = -0.475341651812247
$ws2.Cells.Item(30, 127).Value2 = 1.64778578784757
$ws2.Cells.Item(30, 128).Value2 = -1.06044538706257
$ws2.Cells.Item(32, 128).Value2 = -7.49279538904899
$ws2.Cells.Item(33, 127).Value2 = 2.95790671217292
$ws2.Cells.Item(33, 128).Value2 = 1.94300518134717
$ws2.Cells.Item(34, 128).Value2 = 5.68862275449101
$ws2.Cells.Item(35, 127).Value2 = 0.25641025641026
$ws2.Cells.Item(35, 128).Value2 = -0.00000000000002199822711331579994671985
$ws2.Cells.Item(36, 127).Value2 = 0.729927007299285
$ws2.Cells.Item(36, 128).Value2 = 3.66935483870969
$ws2.Cells.Item(37, 127).Value2 = 11.6352201257862
$ws2.Cells.Item(37, 128).Value2 = 7.75577557755775
$ws2.Cells.Item(38, 127).Value2 = 3.33594361785436
$ws2.Cells.Item(38, 128).Value2 = 1.7961807525052
$ws2.Cells.Item(39, 128).Value2 = 0.721561969439736
$ws2.Cells.Item(40, 127).Value2 = -2.7027027027027
$ws2.Cells.Item(40, 128).Value2 = -6.4516129032258
$ws2.Cells.Item(41, 127).Value2 = 0.562271577171774
$ws2.Cells.Item(41, 128).Value2 = 1.39494333042719
$ws2.Cells.Item(42, 127).Value2 = 1.73584905660378
$ws2.Cells.Item(42, 128).Value2 = 2.61437908496731
$ws2.Cells.Item(43, 127).Value2 = 7.19152157456472
$ws2.Cells.Item(43, 128).Value2 = 2.39190432382705
$ws2.Cells.Item(44, 127).Value2 = 1.07023411371237
$ws2.Cells.Item(44, 128).Value2 = 0.930059523809524
$ws2.Cells.Item(45, 127).Value2 = 1.48698884758365
$ws2.Cells.Item(45, 128).Value2 = -2.13675213675214
$ws2.Cells.Item(46, 127).Value2 = 0.257400257400261
$ws2.Cells.Item(46, 128).Value2 = 0
$ws2.Cells.Item(47, 127).Value2 = 2.10210210210211
$ws2.Cells.Item(47, 128).Value2 = 1.06761565836298
$ws2.Cells.Item(48, 127).Value2 = 0.626304801670155
$ws2.Cells.Item(48, 128).Value2 = 0.899820035992801
$ws2.Cells.Item(49, 127).Value2 = 2.31046931407941
$ws2.Cells.Item(49, 128).Value2 = 1.50196115947575
$ws2.Cells.Item(50, 128).Value2 = 6.01577909270216
$ws2.Cells.Item(51, 127).Value2 = 8.08823529411764
$ws2.Cells.Item(51, 128).Value2 = 5.24193548387099
$ws2.Cells.Item(52, 127).Value2 = 2.27424749163878
$ws2.Cells.Item(52, 128).Value2 = 3.09433962264151
$ws2.Cells.Item(53, 127).Value2 = 3.35542667771334
$ws2.Cells.Item(53, 128).Value2 = 0.43630017452007
$ws2.Cells.Item(54, 127).Value2 = -0.356506238859172
$ws2.Cells.Item(54, 128).Value2 = 1.57170923379176
$ws2.Cells.Item(55, 127).Value2 = 1.77038626609441
$ws2.Cells.Item(55, 128).Value2 = -1.05944673337256
$ws2.Cells.Item(56, 127).Value2 = 2.93040293040293
$ws2.Cells.Item(56, 128).Value2 = 0

